$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Collapse stray double-spaces (after periods) introduced by prior edits.
Replace-Text "corruption issues.  Users" "corruption issues. Users"
Replace-Text "man-in-the-middle attacks.  Also" "man-in-the-middle attacks. Also"
Replace-Text "no correspondence.  When users" "no correspondence. When users"
Replace-Text "to the public.  We hope" "to the public. We hope"
Replace-Text "with information.  Therefore" "with information. Therefore"
Replace-Text "corrupt directly.  Instead" "corrupt directly. Instead"
Replace-Text "data transparency.  Dossants" "data transparency. Dossants"
Replace-Text "publish their findings.  It is" "publish their findings. It is"
Replace-Text "improve our platform.  Dossa welcomes" "improve our platform. Dossa welcomes"

# Update the bug tracker URL.
Replace-Text "https://github.com/dossa/tickets" "https://github.com/thachp/dossa/issues"
